$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: update indices for columns B:E
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2: updated data values for columns B:E
$ws.Range("B2").Value = 49.301159418377814
$ws.Range("C2").Value = -6.6144264138497819
$ws.Range("D2").Value = 0.59876391186929201
$ws.Range("E2").Value = 8.8685582203942523

# Row 3: updated data values for columns B:E
$ws.Range("B3").Value = 46.018332723408086
$ws.Range("C3").Value = 9.1606185307708188
$ws.Range("D3").Value = -19.797397289708901
$ws.Range("E3").Value = 18.036638296737841

# Update the selection to match the new range
$ws.Range("B1:E3").Select()
